$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Range("F3").Value = 413
$ws1.Range("F5").Value = 58
$ws1.Range("F6").Value = 3712
$ws1.Range("F8").Value = 2499
$ws1.Range("F9").Value = 60
$ws1.Range("F10").Value = 2971
$ws1.Range("F13").Value = 2262
$ws1.Range("F17").Value = 422
$ws1.Range("F21").Value = 291
$ws1.Range("F22").Value = 308
$ws1.Range("F28").Value = 142
$ws1.Range("F30").Value = 4094
$ws1.Range("F31").Value = 3666
$ws1.Range("F32").Value = 56
$ws1.Range("F35").Value = 442
$ws1.Range("F38").Value = 138
$ws2.Range("F4").Value = 185
$ws2.Range("F16").Value = 193
$ws3.Range("F2").Value = 1013
$ws3.Range("F4").Value = 2197
$ws4.Range("F3").Value = 1013
$ws4.Range("F5").Value = 413
$ws4.Range("F8").Value = 185
$ws4.Range("F9").Value = 58
$ws4.Range("F11").Value = 3712
$ws4.Range("F13").Value = 2499
$ws4.Range("F14").Value = 60
$ws4.Range("F15").Value = 2971
$ws4.Range("F17").Value = 2262
$ws4.Range("F21").Value = 422
$ws4.Range("B23").Value = '2024-09-15'
$ws4.Range("C23").Value = '北京· EXA·第二届帝都百合only'
$ws4.Range("D23").Value = '永外高庄138号 大红门国际会展中心'
$ws4.Range("E23").Value = '2024.09.15 09:30-09.15 16:00'
$ws4.Range("F23").Value = 332
$ws4.Range("G23").Value = 68
$ws4.Range("H23").Value = 'https://show.bilibili.com/platform/detail.html?id=86477'
$ws4.Range("I23").Value = '//i2.hdslb.com/bfs/openplatform/202405/LgmTjud21716883789133.jpeg'
$ws4.Range("C24").Value = '北京·ICOS SP漫展04动漫节'
$ws4.Range("D24").Value = '金蝉西路甲1号（地铁七号线南楼梓庄站） 北京酷车国际汇展中心'
$ws4.Range("E24").Value = '2024.09.15 09:00-09.16 17:00'
$ws4.Range("F24").Value = 308
$ws4.Range("G24").Value = 80
$ws4.Range("H24").Value = 'https://show.bilibili.com/platform/detail.html?id=90286'
$ws4.Range("I24").Value = '//i2.hdslb.com/bfs/openplatform/202408/tPazRaBV1722595834650.jpeg'
$ws4.Range("C25").Value = '北京·MQ&THEBONE首届怀旧同人only'
$ws4.Range("D25").Value = '安定路5号院(安贞门地铁站A西北口步行420米) 北京北投购物公园'
$ws4.Range("E25").Value = '2024.09.15 10:00-09.16 17:00'
$ws4.Range("F25").Value = 631
$ws4.Range("G25").Value = 6.6
$ws4.Range("H25").Value = 'https://show.bilibili.com/platform/detail.html?id=90096'
$ws4.Range("I25").Value = '//i0.hdslb.com/bfs/openplatform/202407/BBmePlWJ1722409048209.jpeg'
$ws4.Range("C26").Value = '北京·原神only4.0同人展'
$ws4.Range("D26").Value = '北花园路1号 超级蜂巢'
$ws4.Range("E26").Value = '2024.09.15 10:00-09.15 17:00'
$ws4.Range("F26").Value = 1372
$ws4.Range("G26").Value = 68
$ws4.Range("H26").Value = 'https://show.bilibili.com/platform/detail.html?id=87564'
$ws4.Range("I26").Value = '//i0.hdslb.com/bfs/openplatform/202407/EfEAeJDS1720776874376.jpeg'
$ws4.Range("B27").Value = '2024-09-16'
$ws4.Range("C27").Value = '北京·AINI二次元派对【免票展会】'
$ws4.Range("D27").Value = '新村街道丰科路6号F1-102-103 万达广场(丰科店)'
$ws4.Range("E27").Value = '2024.09.16 10:00-09.16 18:00'
$ws4.Range("F27").Value = 34
$ws4.Range("G27").Value = 50
$ws4.Range("H27").Value = 'https://show.bilibili.com/platform/detail.html?id=90730'
$ws4.Range("I27").Value = '//i0.hdslb.com/bfs/openplatform/202408/9SUINRO61723558972754.jpeg'
$ws4.Range("C28").Value = '北京·原神×星穹铁道only2.0同人展'
$ws4.Range("D28").Value = '高碑店东路超级蜂巢 5G直播基地'
$ws4.Range("E28").Value = '2024.09.16 10:00-09.16 17:00'
$ws4.Range("F28").Value = 1283
$ws4.Range("G28").Value = 68
$ws4.Range("H28").Value = 'https://show.bilibili.com/platform/detail.html?id=88285'
$ws4.Range("I28").Value = '//i1.hdslb.com/bfs/openplatform/202406/iWlE3Q9X1719554169582.jpeg'
$ws4.Range("B29").Value = '2024-09-17'
$ws4.Range("C29").Value = '北京·双男主only之皎皎秋月夜'
$ws4.Range("D29").Value = '太平庄中街西端 北京天通苑黄河京都会议中心'
$ws4.Range("E29").Value = '2024.09.17 10:00-09.17 17:00'
$ws4.Range("F29").Value = 117
$ws4.Range("G29").Value = 79
$ws4.Range("H29").Value = 'https://show.bilibili.com/platform/detail.html?id=89763'
$ws4.Range("I29").Value = '//i1.hdslb.com/bfs/openplatform/202407/nUiFpHBb1721723099117.jpeg'
$ws4.Range("F30").Value = 142
$ws4.Range("B32").Value = '2024-09-22'
$ws4.Range("C32").Value = '北京·次元音浪Million Live⏤番音集结'
$ws4.Range("D32").Value = '学清路38号金码大厦B座 北京想象空间'
$ws4.Range("E32").Value = '2024.09.22 13:00-09.22 16:00'
$ws4.Range("F32").Value = 16
$ws4.Range("G32").Value = 88
$ws4.Range("H32").Value = 'https://show.bilibili.com/platform/detail.html?id=90657'
$ws4.Range("I32").Value = '//i1.hdslb.com/bfs/openplatform/202408/Fn9CSOmf1723477511986.jpeg'
$ws4.Range("F33").Value = 4094
$ws4.Range("F34").Value = 3666
$ws4.Range("F35").Value = 56
$ws4.Range("B36").Value = '2024-10-10'
$ws4.Range("C36").Value = '北京·黑白键上的音乐地图——孩子们的钢琴协奏曲之夜'
$ws4.Range("D36").Value = '北新华街1号 北京音乐厅'
$ws4.Range("E36").Value = '2024.10.10 19:30-10.10 21:00'
$ws4.Range("F36").Value = 1
$ws4.Range("G36").Value = 153
$ws4.Range("H36").Value = 'https://show.bilibili.com/platform/detail.html?id=86881'
$ws4.Range("I36").Value = '//i1.hdslb.com/bfs/openplatform/202406/K3oihoH91717474488019.jpeg'
$ws4.Range("B37").Value = '2024-10-11'
$ws4.Range("C37").Value = '北京·官方唯一授权-周杰伦作品国风音乐会'
$ws4.Range("D37").Value = '西直门外大街135号  北展剧场'
$ws4.Range("E37").Value = '2024.10.11 19:30-10.11 21:00'
$ws4.Range("F37").Value = 14
$ws4.Range("G37").Value = 126
$ws4.Range("H37").Value = 'https://show.bilibili.com/platform/detail.html?id=88666'
$ws4.Range("I37").Value = '//i1.hdslb.com/bfs/openplatform/202407/2KgWinEn1720077808243.jpeg'
$ws4.Range("F38").Value = 442
$ws4.Range("F44").Value = 138
$ws4.Range("F49").Value = 193
